$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell J1 = "Q8" (matching the style of the other header cells)
$ws.Range("J1").Value = "Q8"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Update numeric error values across the data grid (B2:J15)
$ws.Range("B2").Value = -0.06845131779842006
$ws.Range("C2").Value = -2.344079640195531
$ws.Range("D2").Value = 0.000005587651372707114
$ws.Range("E2").Value = -0.6233623182505283
$ws.Range("F2").Value = 4.887355777631257
$ws.Range("G2").Value = -0.07746478018816561
$ws.Range("H2").Value = -0.223199686849668
$ws.Range("B3").Value = -1.844079640195545
$ws.Range("C3").Value = 0.5000055876513585
$ws.Range("D3").Value = -0.1233623182505426
$ws.Range("E3").Value = 5.387355777631242
$ws.Range("F3").Value = 0.4225352198118202
$ws.Range("G3").Value = 0.2768003131503178
$ws.Range("B4").Value = 0.7000055876513613
$ws.Range("C4").Value = 0.07663768174946028
$ws.Range("D4").Value = 5.587355777631245
$ws.Range("E4").Value = 0.622535219811823
$ws.Range("F4").Value = 0.4768003131503207
$ws.Range("G4").Value = 0.3526870962688803
$ws.Range("H4").Value = 1.126878279600263
$ws.Range("I4").Value = -0.334675296103242
$ws.Range("J4").Value = 0.265633061926664
$ws.Range("B5").Value = -1.623362318250543
$ws.Range("C5").Value = 3.887355777631242
$ws.Range("D5").Value = -1.07746478018818
$ws.Range("E5").Value = -1.223199686849682
$ws.Range("F5").Value = -1.347312903731122
$ws.Range("G5").Value = -0.5731217203997403
$ws.Range("H5").Value = -2.034675296103245
$ws.Range("I5").Value = -1.434366938073339
$ws.Range("B6").Value = 4.387355777631228
$ws.Range("C6").Value = -0.577464780188194
$ws.Range("D6").Value = -0.7231996868496964
$ws.Range("E6").Value = -0.8473129037311367
$ws.Range("F6").Value = -0.07312172039975451
$ws.Range("G6").Value = -1.534675296103259
$ws.Range("H6").Value = -0.934366938073353
$ws.Range("B7").Value = 2.722535219811803
$ws.Range("C7").Value = 2.576800313150301
$ws.Range("D7").Value = 2.452687096268861
$ws.Range("E7").Value = 3.226878279600243
$ws.Range("F7").Value = 1.765324703896738
$ws.Range("G7").Value = 2.365633061926644
$ws.Range("B8").Value = -0.2231996868496964
$ws.Range("C8").Value = -0.3473129037311367
$ws.Range("D8").Value = 0.4268782796002455
$ws.Range("E8").Value = -1.034675296103259
$ws.Range("F8").Value = -0.434366938073353
$ws.Range("G8").Value = -0.41600913674678
$ws.Range("H8").Value = 0.3423454266220887
$ws.Range("I8").Value = -0.7919400257838731
$ws.Range("B9").Value = -0.1473129037311339
$ws.Range("C9").Value = 0.6268782796002483
$ws.Range("D9").Value = -0.8346752961032562
$ws.Range("E9").Value = -0.2343669380733502
$ws.Range("F9").Value = -0.2160091367467772
$ws.Range("G9").Value = 0.5423454266220915
$ws.Range("H9").Value = -0.5919400257838703
$ws.Range("B10").Value = 0.3268782796002512
$ws.Range("C10").Value = -1.134675296103253
$ws.Range("D10").Value = -0.5343669380733473
$ws.Range("E10").Value = -0.5160091367467743
$ws.Range("F10").Value = 0.2423454266220944
$ws.Range("G10").Value = -0.8919400257838674
$ws.Range("B11").Value = -0.8346752961032562
$ws.Range("C11").Value = -0.2343669380733502
$ws.Range("D11").Value = -0.2160091367467772
$ws.Range("E11").Value = 0.5423454266220915
$ws.Range("F11").Value = -0.5919400257838703
$ws.Range("B12").Value = -0.3343669380733587
$ws.Range("C12").Value = -0.3160091367467857
$ws.Range("D12").Value = 0.442345426622083
$ws.Range("E12").Value = -0.6919400257838788
$ws.Range("B13").Value = -0.3160091367467857
$ws.Range("C13").Value = 0.442345426622083
$ws.Range("D13").Value = -0.6919400257838788
$ws.Range("B14").Value = 0.3423454266220887
$ws.Range("C14").Value = -0.7919400257838731
$ws.Range("B15").Value = -0.2919400257838873

Write-Output "Edit complete"
